$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text columns to be treated as text (avoid date/number auto-conversion)
$ws.Range("B2").NumberFormat = "@"

$ws.Range("A2").Value = "tA6HgoO8"
$ws.Range("B2").Value = "10/11/2024"
$ws.Range("C2").Value = "23:00"
$ws.Range("D2").Value = "MEXICO - LIGA DE EXPANSION MX"
$ws.Range("E2").Value = "Tapatio"
$ws.Range("F2").Value = "Tepatitlan de Morelos"

# Restore default style on B2 after forcing text format
$ws.Range("B2").Style = "Normal"

$ws.Range("G2").Value = 1.87
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 2.42
$ws.Range("K2").Value = 2.1
$ws.Range("L2").Value = 4.35
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 7.5
$ws.Range("O2").Value = 1.35
$ws.Range("P2").Value = 2.7
$ws.Range("Q2").Value = 2.02
$ws.Range("R2").Value = 1.62
$ws.Range("S2").Value = 1.39
$ws.Range("T2").Value = 2.55
$ws.Range("U2").Value = 1.87
$ws.Range("V2").Value = 1.75
$ws.Range("W2").Value = 6.2
$ws.Range("X2").Value = 8.25
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 15.5
$ws.Range("AA2").Value = 16.5
$ws.Range("AB2").Value = 32
$ws.Range("AC2").Value = 8.5
$ws.Range("AD2").Value = 6.4
$ws.Range("AE2").Value = 16.5
$ws.Range("AF2").Value = 90
$ws.Range("AG2").Value = 800
$ws.Range("AH2").Value = 10.25
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 13.5
$ws.Range("AK2").Value = 65
$ws.Range("AL2").Value = 40
$ws.Range("AM2").Value = 50
$ws.Range("AN2").Value = 3.65
$ws.Range("AO2").Value = 9.25
$ws.Range("AP2").Value = 18.5
$ws.Range("AQ2").Value = 32
$ws.Range("AR2").Value = 70
$ws.Range("AS2").Value = 250
$ws.Range("AT2").Value = 2.5
$ws.Range("AU2").Value = 7.2
$ws.Range("AV2").Value = 65
$ws.Range("AW2").Value = 5.7
$ws.Range("AX2").Value = 22
$ws.Range("AY2").Value = 28
$ws.Range("AZ2").Value = 120
$ws.Range("BA2").Value = 150
$ws.Range("BB2").Value = 350
$ws.Range("BC2").Value = 51
$ws.Range("BD2").Value = 51

# Remove the now-obsolete third row (original row 3 data is now in row 2)
$ws.Rows.Item(3).Delete()
